# Copy changes to template: rename the "Gender" header to "Sex".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The header row (row 7) has "Gender" in column H; change it to "Sex".
$ws.Range("H7").Value = "Sex"

# Leave the selection on the edited cell, matching the authored change.
$ws.Range("H7").Select()
